# Stage 1: update companies data
#
# The data rows (3-28) of the "relevant companies" sheet are being
# re-ordered: each destination row receives the full record (columns
# A-K) that used to live at a different source row. Capture every row's
# current values first (so source data isn't clobbered before it is
# read), then write the values back out in their new positions.
#
# All columns in this sheet are plain text (company numbers, SIC codes,
# dates, times, etc. are all stored as text, not numbers/dates), so
# number/date-looking values are re-entered with a leading apostrophe to
# keep Excel from auto-converting them to a numeric or date type.

function Looks-NumericOrDate($s) {
    if ($null -eq $s -or $s -eq "") { return $false }
    if ($s -match '^-?[0-9]+(\.[0-9]+)?$') { return $true }
    if ($s -match '^\d{1,4}[/-]\d{1,2}[/-]\d{1,4}$') { return $true }
    return $false
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# destination row -> source row (data that should end up there)
$mapping = @{
    3  = 20
    4  = 21
    5  = 19
    6  = 18
    7  = 16
    8  = 15
    9  = 17
    10 = 26
    11 = 25
    12 = 24
    13 = 28
    14 = 27
    15 = 22
    16 = 23
    17 = 7
    18 = 8
    19 = 4
    20 = 3
    21 = 6
    22 = 5
    23 = 12
    24 = 14
    25 = 13
    26 = 9
    27 = 11
    28 = 10
}

$firstCol = 1   # A
$lastCol  = 11  # K
$skipCols = @(5)  # column E ("Source") is blank in every row - leave it alone

# Snapshot every source row's text before any writes happen.
$snapshot = @{}
for ($r = 3; $r -le 28; $r++) {
    $rowData = @{}
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowData[$c] = $ws.Cells.Item($r, $c).Text
    }
    $snapshot[$r] = $rowData
}

# Write the snapshotted values back into their new (destination) rows.
foreach ($destRow in $mapping.Keys) {
    $srcRow = $mapping[$destRow]
    $rowData = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        if ($skipCols -contains $c) { continue }
        $val = $rowData[$c]
        $cell = $ws.Cells.Item($destRow, $c)
        if (Looks-NumericOrDate $val) {
            # Force text interpretation (like a leading apostrophe in the
            # UI) so company numbers / dates / SIC codes stay text instead
            # of being auto-converted to a number or date serial, then
            # drop back to the default "Normal" style so no stray
            # quote-prefix formatting is left behind on the cell.
            $cell.Value = "'" + $val
            $cell.Style = "Normal"
        } else {
            $cell.Value = $val
        }
    }
}
